$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-24 Monday" "2025-11-25 Tuesday"

Replace-Text "18÷9=" "92÷6="
Replace-Text "61÷5=" "80÷6="
Replace-Text "13÷8=" "37÷2="
Replace-Text "27÷3=" "59÷3="
Replace-Text "43÷2=" "43÷8="

Replace-Text "15÷9=" "97÷9="
Replace-Text "37÷9=" "15÷6="
Replace-Text "11÷6=" "20÷5="
Replace-Text "25÷6=" "66÷5="
Replace-Text "89÷4=" "14÷3="

Replace-Text "82÷7=" "47÷6="
Replace-Text "53÷5=" "36÷2="
Replace-Text "19÷7=" "67÷6="
Replace-Text "32÷2=" "99÷5="
Replace-Text "36÷8=" "65÷7="

Replace-Text "14÷5=" "52÷7="
Replace-Text "29÷8=" "79÷5="
Replace-Text "60÷6=" "43÷8="
Replace-Text "23÷7=" "53÷4="
Replace-Text "99÷9=" "36÷8="

Replace-Text "61÷2=" "97÷8="
Replace-Text "86÷3=" "40÷3="
Replace-Text "60÷8=" "49÷7="
Replace-Text "81÷7=" "64÷5="
Replace-Text "10÷2=" "95÷4="
